$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.428.02'
$ws.Range("E2").Value = '  +1.12%  '

# Row 3
$ws.Range("D3").Value = '2.592.24'
$ws.Range("E3").Value = '  -0.04%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.72'
$ws.Range("E5").Value = '  +3.54%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.16'
$ws.Range("E6").Value = '  +1.01%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("E8").Value = '  +0.27%  '

# Row 9
$ws.Range("D9").Value = '2.603.25'
$ws.Range("E9").Value = '  +0.00%  '

# Row 10
$ws.Range("E10").Value = '  -1.49%  '

# Row 11
$ws.Range("E11").Value = '  +3.91%  '

# Row 12
$ws.Range("E12").Value = '  +11.24%  '

# Row 13
$ws.Range("E13").Value = '  +3.62%  '

# Row 14
$ws.Range("D14").Value = '3.047.89'
$ws.Range("E14").Value = '  -0.16%  '

# Row 15
$ws.Range("D15").Value = '59.445.71'
$ws.Range("E15").Value = '  +1.23%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.62'
$ws.Range("E16").Value = '  +8.82%  '

# Row 17
$ws.Range("E17").Value = '  +4.68%  '

# Row 18
$ws.Range("D18").Value = '2.598.72'
$ws.Range("E18").Value = '  -0.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.55'
$ws.Range("E19").Value = '  +1.82%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '338.24'
$ws.Range("E20").Value = '  +0.45%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.28'
$ws.Range("E21").Value = '  +2.31%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.22'
$ws.Range("E22").Value = '  +0.85%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.11%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.60'
$ws.Range("E24").Value = '  -3.28%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.456'
$ws.Range("E25").Value = '  +6.62%  '

# Row 26
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.162'
$ws.Range("E26").Value = '  +2.13%  '

# Row 27
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.30'
$ws.Range("E28").Value = '  +2.55%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0786'
$ws.Range("E29").Value = '  +4.69%  '

# Row 30
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$ws.Range("E31").Value = '  +1.13%  '

# Row 32
$ws.Range("E32").Value = '  +1.75%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '159.44'
$ws.Range("E33").Value = '  +3.36%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.11'
$ws.Range("E34").Value = '  +0.99%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.05'
$ws.Range("E35").Value = '  +3.40%  '

# Row 36
$ws.Range("E36").Value = '  +2.29%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.887'
$ws.Range("E37").Value = '  +5.15%  '

# Row 38
$ws.Range("E38").Value = '  -1.26%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.50'
$ws.Range("E39").Value = '  +2.85%  '

# Row 40
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.13'
$ws.Range("E40").Value = '  +0.80%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '295.80'
$ws.Range("E41").Value = '  +4.51%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.69'
$ws.Range("E42").Value = '  +2.74%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0980'
$ws.Range("E44").Value = '  +2.68%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.598'
$ws.Range("E45").Value = '  -0.04%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0540'
$ws.Range("E46").Value = '  +1.48%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.32'
$ws.Range("E47").Value = '  +3.54%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.63'
$ws.Range("E48").Value = '  -0.01%  '

# Row 49
$ws.Range("E49").Value = '  +2.86%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.20'
$ws.Range("E50").Value = '  +5.75%  '

# Row 51
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.952.17'
$ws.Range("E51").Value = '  +0.55%  '
